$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.714.68'
$ws.Range("E2").Value = '  -1.05%  '
$ws.Range("D3").Value = '3.496.85'
$ws.Range("E3").Value = '  -1.56%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '595.63'
$ws.Range("E5").Value = '  -1.62%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '142.43'
$ws.Range("E6").Value = '  -1.49%  '
$ws.Range("D7").Value = '3.496.77'
$ws.Range("E7").Value = '  -1.53%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("E8").Value = '  -0.27%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.498'
$ws.Range("E9").Value = '  +0.33%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.62'
$ws.Range("E11").Value = '  -3.94%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.402'
$ws.Range("E12").Value = '  -2.69%  '
$ws.Range("D13").Value = '4.090.33'
$ws.Range("E13").Value = '  -1.55%  '
$ws.Range("E14").Value = '  -4.10%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '28.53'
$ws.Range("E15").Value = '  -4.89%  '
$ws.Range("D16").Value = '3.486.71'
$ws.Range("E16").Value = '  -2.11%  '
$ws.Range("E17").Value = '  +1.19%  '
$ws.Range("D18").Value = '65.717.75'
$ws.Range("E18").Value = '  -1.12%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.90'
$ws.Range("E19").Value = '  -5.95%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.14'
$ws.Range("E20").Value = '  -0.43%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.27'
$ws.Range("E21").Value = '  -3.79%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '410.94'
$ws.Range("E22").Value = '  -4.64%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.590'
$ws.Range("E23").Value = '  -3.13%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '77.46'
$ws.Range("D25").Value = '3.636.58'
$ws.Range("E25").Value = '  -1.50%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  +0.05%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000114'
$ws.Range("E27").Value = '  -4.75%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.96'
$ws.Range("E28").Value = '  -2.80%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.42'
$ws.Range("E29").Value = '  -3.35%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.63'
$ws.Range("E30").Value = '  -4.29%  '
$ws.Range("E31").Value = '  +0.03%  '
$ws.Range("D32").Value = '3.494.85'
$ws.Range("E32").Value = '  -1.47%  '
$ws.Range("E33").Value = '  -0.60%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '24.10'
$ws.Range("E34").Value = '  -4.86%  '
$ws.Range("E35").Value = '  -0.01%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '7.41'
$ws.Range("E36").Value = '  -5.64%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '174.70'
$ws.Range("E37").Value = '  -0.01%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.23'
$ws.Range("E38").Value = '  -15.22%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.16'
$ws.Range("E39").Value = '  -7.51%  '
$ws.Range("E40").Value = '  -9.39%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0812'
$ws.Range("E41").Value = '  -4.25%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.00'
$ws.Range("E42").Value = '  -3.59%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.845'
$ws.Range("E43").Value = '  -4.76%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '45.19'
$ws.Range("E44").Value = '  -2.00%  '
$ws.Range("E45").Value = '  -8.61%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.999'
$ws.Range("E46").Value = '  +0.02%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.38'
$ws.Range("E47").Value = '  -6.49%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.03'
$ws.Range("E48").Value = '  -1.56%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '22.24'
$ws.Range("E49").Value = '  -5.13%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.07'
$ws.Range("E50").Value = '  -9.14%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '22.62'
$ws.Range("E51").Value = '  -9.72%  '
